# Edit Review_423.docx per the target diff:
#  - Update the date in the title line (20.03.25 -> 19.03.25)
#  - Update the paper title (softmax... -> EFFICIENTLY LEARNING AT TEST-TIME...)
#  - Rewrite the body paragraphs (new paper summary about SIFT / active fine-tuning)
#  - Append a dozen new paragraphs for the rest of the new review
#  - Update the trailing arXiv link

$d = $word.ActiveDocument

# --- 1. Title paragraph: two runs of text split by a <w:br/> -----------------
$d.Content.Find.Execute(
    'המאמר היומי של מייק - 20.03.25', $true, $false, $false, $false, $false,
    $true, 1, $false, 'המאמר היומי של מייק - 19.03.25', 2) | Out-Null

$d.Content.Find.Execute(
    'softmax is not enough (for sharp out-of-distribution)', $true, $false, $false, $false, $false,
    $true, 1, $false, 'EFFICIENTLY LEARNING AT TEST-TIME: ACTIVE FINE-TUNING OF LLMS', 2) | Out-Null

# --- 2. Paragraph 2 (index 1) --------------------------------------------------
$old2 = 'המאמר הזה מציעה שיטה לשיפור ביצועי ההכללה עבור מודלי טרנספורמרים מזווית די לא צפויה. המחברים מציעים שיטה להתמודדות עם מה שנקרא דיספרסיה (או פיזור בעברית) של מקדמים ה-attention בטרנספורמרים. זה מתבטא למשל באי יכולת (לפי המאמר) של הטרנספורמרים למקד את מקדמי ה-attention במספר טוקנים קטן (יחסית לאורך הסדרה). זה חשוב למשל בשאלות כמו מציאת מקסימומים של סדרת מספרים נתונה או שאלות בסגנון ״מחט בערימת השחת״ (needle in a haystack) כאשר המודל מתבקש מקטע קצר לא קשור בטקסט מסוים (יחסית ארוך).'
$new2 = 'בתקופה האחרונה השיטה הכי פופולרית להתאמת מודלי שפה למשימה ספציפית היא למידה in-context או ICL. בגדול אנו מספקים למודל, בתוך הפרומפט, כמה דוגמאות לביצוע משימה והמודל ״לומד״ איך לבצע אותה ללא שום שינוי במשקליו. ICL מתאפשר עקב האופי האדפטיבי של הטרנספורמרים (מנגנון ה-attention בתוכו) המצליחים ״לעדכן את אופן החישוב שלו״ כפונקציה של קלט. '
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- 3. Paragraph 3 (index 2) --------------------------------------------------
$old3 = 'המחברים טוענים שאחת הסיבות לבעיות אלו היא פיזור מקדמי ה-attention במנגנון הטרנספורמרים. מקדמים אלו מחושבים עם פונקצית סופטמקס ה״מנרמלת״ את המכפלות הפנימיות של וקטורי K ו-Q עבור כל טוקני הסדרה. לפי המאמר הבעיה קשורה לכך שעבור קונטקסטים ארוכים לסופטמקס במיוחד בטרנספורמרים העמוקים יש ״נטיה למרוח את פלט הסופטמקס״. '
$new3 = 'המאמר דן בשיטה אחרת לאדפטציה של מודל למשימה נתונה בזמן טסט(המאמר קצת מערבב את המושג של  טסט ואינפרנס) המערב fine-tune קליל של המודל על סמך הפרומפט שמוזן אליו. להבדיל מ-ICL השיטה המוצעת (SIFT(Selects Informative data for Fine-Tuning כן משנה את משקלי המודל (מבצעת צעד אחד של מורד הגרדיאנט - gradient descent). למעשה SIFT (ד״א יש שיטה בשם כזה גם בעיבוד תמונה מהעידן לפני הרשתות) מציעה שיטה לבחירה של דוגמאות מהדאטהסט לפיין טיון של מודל עבור פרומפט נתון.'
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- 4. Paragraph 4 (index 3) --------------------------------------------------
$old4 = 'אחת הדרכים להתמודד עם התופעה הזו היא להוריד את הטמפרטורה אבל זה עלול להעלות סיכוי לשגיאה במקרים בהם הלוגיט (משקל attention לא מנורמל) של הטוקן הנכון יותר קטן מהלוגית המקסימלי. כדי להתמודד עם התופעה המבחרים הציעו גרסה חדשה של סופטמקס בה הטמפרטורה תלויה באנטרופיה של הטוקנים. '
$new4 = 'המחברים טוענים שבחירת דוגמאות הכי קרובות לפרומפט במרחב הלטנטי מבחינת מרחק קוסיין או מכפלה פנימית(nearest neighbors or NN) היא תת-אופטימלית ועלולה להביא דוגמאות מיותרות הפוגעות בביצועי פיין טיון. במקום לשלוף דוגמאות הדומות ביותר לפרומפט, SIFT בוחרת את אלו שמספקות את מירב המידע החדש, וכך משיגה התאמה טובה יותר של המודל עם מינימום חישובים נוספים.'
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# --- 5. Paragraph 5 (index 4) --------------------------------------------------
$old5 = 'הם אימנו מודל עבור מקרים שבהם הלוגיט של הטוקן הנכון אינו מקסימלי כאשר המטרה היתה למקסם את הסתברות הדגימה של הטוקן הנכון (אחרי מנגנון ה-attention ו-FFN). מטרת המודל היתה לחשב ערך אופטימלי של טמפרטורה כפונקצייה של אנטרופיית של משקלי attention לא מנורמלים. הנוסחה של הטמפרטורה יצאה הופכית (1 חלקי) של פולימום מחזקה 4. אציין כי הטמפרטורה מחושבת בזמן האינפרנס כתלות באינטרופיית הטוקנים לפי המודל הזה.'
$new5 = 'הגישה המוצעת מערבת שיעורך רמת אי ודאות של תשובת המודל בהינתן הדוגמאות שבחרנו ל-FT (לאחר FT הכוונה). בפרק הבא אסביר למה זה חשוב בעצם.'
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# --- 6. Paragraph 6 (index 5) --------------------------------------------------
$old6 = 'המחברים הראו אמפירית כי עם הטמרטורה האדפטיבית מקטינה פיזור משקלי ה-attention. למרות שהטמפרטורה האדפטיבית האופטימלית יורדת עם עלייה באנטרופיית הלוגיטים היא גורמת לפחות שגיאות של המודל יחסית למקרה שהיא נקבעת באופן קשיח. '
$new6 = 'הערכת אי-וודאות להנחיית FT ולמה זה בכלל חשוב כאן?'
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# --- 7. Paragraph 7 (index 6, was the arXiv link) ------------------------------
$old7 = 'https://arxiv.org/abs/2410.01104'
$new7 = 'שיטות FT רבות מסתמכות על שליפת דוגמאות דומות בהתבסס על דמיון קוסיין או מרחק אוקלידי. אך גישה זו לוקה בחסר: היא אינה מבדילה בין דאטה רלוונטי לזה שמיותר. שתי דוגמאות דומות מאוד עשויות להכיל את אותו מידע, ולכן אחת מהן אינה תורמת לתוצאת FT. כדי לפתור זאת, המחברים מציעים שיטה להערכת אי-הוודאות של המודל בתשובתו לאחר FT . אם המודל בטוח מאוד בתשובתו אחרי FT, הוספת דוגמא לא תשפיע משמעותית. אך אם אי-הוודאות גבוהה, בחירה חכמה של דוגמאות יכולה לשפר את ביצועי המודל משמעותית ןהאתגר הוא למצוא את הדוגמאות הללו ביעילות.'
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# --- 8. Append the new trailing paragraphs -------------------------------------
$newParagraphs = @(
    'מדידת דמיון במרחב הסמוי בעזרת פונקציית קרנל',
    'כאמור הבסיס לשיטת הבחירה של SIFT הוא מדידת הדמיון בין דוגמאות במרחב לטנטי. כדי לכמת את הדמיון הזה, המחברים משתמשים בפונקציית קרנל - שהיא מוגדרת בתור מכפלה פנימית בין הייצוגים הלטנטיים של הדוגמאות. פונקציה זו מקבלת שני רצפים ומחזירה ציון דמיון—גבוה עבור סדרות דומות ונמוך עבור רצפים שונים. בעזרת פונקצית קרנל זו בונים מטריצה קרנל עבור הדוגמאות שנבחרו ל- FT והפרומפט עצמו. לאחר מכן מגדירי מודל דמה (surrogate model) שמטרו לשערך את ביצועי ה-LLM לאחר FT על הדוגמאות שנבחרו.',
    'באמצעות מודל זה בונים (זה קצת כבד מתמטית) את השיערוך של אי וודאות של המודל אחרי הוספה של דוגמא x מהדאטהסט לסט הדוגמאות שעליהם יתבצע הטיוב. בסופו של דבר בוחרים דוגמא הממזערת את אי ודאות עבור הפרומפט ומוסיפים אותה לסט הדוגמאות זה.',
    'במילים פשוטות הגישה המוצעות מאזנת בין שני שיקולים מנוגדים:',
    'רלוונטיות: הדוגמאות הנבחרות צריכות להיות עדיין רלוונטיות לפרומפט.',
    'גיוון: הדוגמאות אינן אמורות להכיל מידע חופף ומיותר.',
    'במקום לבחור דוגמאות בבת אחת, SIFT בוחר כל דוגמה באופן הדרגתי, תוך שימוש בפונקציית קרנל כדי לקבוע את הערך המוסף שלה.',
    'אם מועמד חדש דומה מדי לדוגמאות שנבחרו בעבר, הוא נדחה, מכיוון שהוא אינו מוסיף מידע חדש.',
    'אם המועמד רלוונטי אך מכיל פרטים חדשים, הוא נבחר כדי להפחית את אי-הוודאות.',
    'אם המועמד אינו קשור לפרומפט כלל, הוא נשאר מחוץ לתהליך.',
    'https://arxiv.org/abs/2410.08020'
)

foreach ($t in $newParagraphs) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $t
}
